$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.799.92"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "2.455.70"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'570.04"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("D6").Value = "'146.68"
$ws.Range("E6").Value = "  +1.94%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D9").Value = "2.455.37"
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("E10").Value = "  +1.99%  "
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "'0.355"
$ws.Range("E13").Value = "  +1.72%  "
$ws.Range("D14").Value = "'26.86"
$ws.Range("E14").Value = "  +2.55%  "
$ws.Range("E15").Value = "  +3.22%  "
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("D17").Value = "62.853.69"
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("D18").Value = "2.450.57"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("D19").Value = "'11.35"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("D20").Value = "'7.27"
$ws.Range("E20").Value = "  +6.48%  "
$ws.Range("D21").Value = "'323.09"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("E23").Value = "  +12.90%  "
$ws.Range("D24").Value = "'0.998"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").Value = "'66.16"
$ws.Range("E25").Value = "  -1.89%  "
$ws.Range("D26").Value = "'619.98"
$ws.Range("E26").Value = "  +11.06%  "
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("E28").Value = "  +9.07%  "
$ws.Range("E29").Value = "  +0.99%  "
$ws.Range("D30").Value = "'0.997"
$ws.Range("E30").Value = "  -0.28%  "
$ws.Range("E31").Value = "  +6.03%  "
$ws.Range("D32").Value = "'8.26"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("E33").Value = "  -3.84%  "
$ws.Range("E34").Value = "  +1.84%  "
$ws.Range("E35").Value = "  +6.41%  "
$ws.Range("E36").Value = "  -1.84%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E39").Value = "  -1.49%  "
$ws.Range("D40").Value = "'18.66"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").Value = "'144.61"
$ws.Range("E41").Value = "  -4.83%  "
$ws.Range("E42").Value = "  -1.33%  "
$ws.Range("D43").Value = "'2.61"
$ws.Range("E43").Value = "  +15.66%  "
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "'147.33"
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("E46").Value = "  +2.18%  "
$ws.Range("D47").Value = "'20.68"
$ws.Range("E47").Value = "  +3.56%  "
$ws.Range("D48").Value = "'0.0536"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("D49").Value = "'0.601"
$ws.Range("E49").Value = "  +0.89%  "
$ws.Range("E50").Value = "  +2.31%  "
$ws.Range("E51").Value = "  -0.60%  "
